$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 118 and 119 (resultado / profit columns) ---
$ws.Range("G118").Value = "Acierto"
$ws.Range("H118").Value = 0.8
$ws.Range("G119").Value = "Acierto"
$ws.Range("H119").Value = 1.63

# --- Append new match rows (auto-updated tracker entries) ---
# Force column B to be treated as plain text so date-like strings
# (e.g. "2025-09-14") are not auto-converted into Excel date serials.
$ws.Range("B124:B138").NumberFormat = "@"

$ws.Range("A124").Value = 14679457
$ws.Range("B124").Value = '2025-09-14'
$ws.Range("C124").Value = 'Matyas Cerny'
$ws.Range("D124").Value = 'Cedrik-Marcel Stebe'
$ws.Range("E124").Value = 'Gana Matyas Cerny'
$ws.Range("F124").Value = 4

$ws.Range("A125").Value = 14679460
$ws.Range("B125").Value = '2025-09-14'
$ws.Range("C125").Value = 'Syl Gaxherri'
$ws.Range("D125").Value = 'Dennis Novak'
$ws.Range("E125").Value = 'Gana Syl Gaxherri'
$ws.Range("F125").Value = 13

$ws.Range("A126").Value = 14679455
$ws.Range("B126").Value = '2025-09-14'
$ws.Range("C126").Value = 'Hynek Barton'
$ws.Range("D126").Value = 'Stijn Slump'
$ws.Range("E126").Value = 'Gana Stijn Slump'
$ws.Range("F126").Value = 6.5

$ws.Range("A127").Value = 14679461
$ws.Range("B127").Value = '2025-09-14'
$ws.Range("C127").Value = 'Zdenek Kolar'
$ws.Range("D127").Value = 'Nikita Mashtakov'
$ws.Range("E127").Value = 'Gana Nikita Mashtakov'
$ws.Range("F127").Value = 7.5

$ws.Range("A128").Value = 14679458
$ws.Range("B128").Value = '2025-09-14'
$ws.Range("C128").Value = 'Alex Huszar'
$ws.Range("D128").Value = 'Neil Oberleitner'
$ws.Range("E128").Value = 'Gana Alex Huszar'
$ws.Range("F128").Value = 12

$ws.Range("A129").Value = 14679456
$ws.Range("B129").Value = '2025-09-14'
$ws.Range("C129").Value = 'Matthew William Donald'
$ws.Range("D129").Value = 'Lukas Pokorny'
$ws.Range("E129").Value = 'Gana Matthew William Donald'
$ws.Range("F129").Value = 2.1

$ws.Range("A130").Value = 14679555
$ws.Range("B130").Value = '2025-09-14'
$ws.Range("C130").Value = 'Alexey Vatutin'
$ws.Range("D130").Value = 'Arthur Reymond'
$ws.Range("E130").Value = 'Gana Arthur Reymond'
$ws.Range("F130").Value = 2.75

$ws.Range("A131").Value = 14679556
$ws.Range("B131").Value = '2025-09-14'
$ws.Range("C131").Value = 'Vadym Ursu'
$ws.Range("D131").Value = 'Luca Potenza'
$ws.Range("E131").Value = 'Gana Vadym Ursu'
$ws.Range("F131").Value = 3

$ws.Range("A132").Value = 14679558
$ws.Range("B132").Value = '2025-09-14'
$ws.Range("C132").Value = 'Michael Geerts'
$ws.Range("D132").Value = 'Mikail Alimli'
$ws.Range("E132").Value = 'Gana Mikail Alimli'
$ws.Range("F132").Value = 10.5

$ws.Range("A133").Value = 14679552
$ws.Range("B133").Value = '2025-09-14'
$ws.Range("C133").Value = 'Leo Raquillet'
$ws.Range("D133").Value = 'Aryan Shah'
$ws.Range("E133").Value = 'Gana Leo Raquillet'
$ws.Range("F133").Value = 3.75

$ws.Range("A134").Value = 14679550
$ws.Range("B134").Value = '2025-09-14'
$ws.Range("C134").Value = 'Louis Tessa'
$ws.Range("D134").Value = 'Karan Singh'
$ws.Range("E134").Value = 'Gana Louis Tessa'
$ws.Range("F134").Value = 3.25

$ws.Range("A135").Value = 14679557
$ws.Range("B135").Value = '2025-09-14'
$ws.Range("C135").Value = 'Rodrigo Pacheco Mendez'
$ws.Range("D135").Value = 'Aidan McHugh'
$ws.Range("E135").Value = 'Gana Aidan McHugh'
$ws.Range("F135").Value = 3.75

$ws.Range("A136").Value = 14680989
$ws.Range("B136").Value = '2025-09-14'
$ws.Range("C136").Value = 'Hernan Casanova'
$ws.Range("D136").Value = 'Juan Sebastian Gomez'
$ws.Range("E136").Value = 'Gana Juan Sebastian Gomez'
$ws.Range("F136").Value = 5.5

$ws.Range("A137").Value = 14680992
$ws.Range("B137").Value = '2025-09-14'
$ws.Range("C137").Value = 'Facundo Bagnis'
$ws.Range("D137").Value = 'Federico Aguilar Cardozo'
$ws.Range("E137").Value = 'Gana Federico Aguilar Cardozo'
$ws.Range("F137").Value = 9

$ws.Range("A138").Value = 14680987
$ws.Range("B138").Value = '2025-09-14'
$ws.Range("C138").Value = 'Ignacio Monzon'
$ws.Range("D138").Value = 'Tomas Martinez'
$ws.Range("E138").Value = 'Gana Tomas Martinez'
$ws.Range("F138").Value = 2.1

# Restore default (unstyled) formatting on column B for the new rows so
# the appended cells match the plain, style-less cells used elsewhere
# in the sheet for this column.
$ws.Range("B124:B138").Style = "Normal"
